$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update unit quantities (Inversiones section) ---
$ws.Range("D55").Value = 10
$ws.Range("D60").Value = 3
$ws.Range("D62").Value = 18
$ws.Range("D63").Value = 10

# --- Re-fill E53:E64 with the Total formula so it becomes one shared formula range ---
$ws.Range("E53:E64").Formula = "=C53*D53"

# --- Add the combined "TOTAL G+I" block below the Inversiones total ---
$ws.Range("E70").Value = "TOTAL G+I"
$ws.Range("E71").Formula = "=SUM(E48,E67)"

# Bold + fill for the "TOTAL" style cells (matches the existing TOTAL rows) and the new TOTAL G+I label
$ws.Range("E47:E48,E66:E67,E70").Font.Bold = $true
$ws.Range("E47:E48,E66:E67,E70,E71").Interior.ThemeColor = 4

# --- Update the view so the new rows are visible/selected (cosmetic) ---
$win = $excel.ActiveWindow
$win.ScrollRow = 52
$win.ScrollColumn = 1
$ws.Range("E71").Select() | Out-Null

Write-Host "done"
